# Generate Report for Handoff
#
# Rows 4-7 (files 4a511410-, 5f0fe1e7-, 9e56768c-, b6a21be4-) were
# re-generated for handoff: their priority flips from "low" to "ht", the
# zh-cn handoff timestamp moves from 20:35:13 to 20:35:37, and the Overview
# sheet's "Latest HO Xliff Generate Date" moves from 20:35:19 to 20:35:41
# (the de-de sheet's "Latest Handoff Datetime" for those same rows shares
# that exact string in the workbook's string table, so it moves too).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 4-7
$overview.Range("G4:G7").Value = "2016-08-30 20:35:41"

# zh-cn sheet: Priority (E) low -> ht, and Latest Handoff Datetime (H)
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-30 20:35:37"

# de-de sheet: Priority (E) low -> ht; Latest Handoff Datetime (H) shares
# the same underlying string as Overview's G4:G7, so it carries the same update
$dede.Range("E4").Value = "ht"
$dede.Range("E5").Value = "ht"
$dede.Range("E6").Value = "ht"
$dede.Range("E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-30 20:35:41"
